# Adds 2023-12-19 YTD violent crime data updates across Citywide, By Neighborhood,
# and individual neighborhood worksheets. Only column J (year 2023) values change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7413
$ws.Range("J3").Value = 7807
$ws.Range("J4").Value = 1697
$ws.Range("J5").Value = 612
$ws.Range("J6").Value = 10651
$ws.Range("J7").Value = 28180

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 466
$ws.Range("J5").Value = 45
$ws.Range("J6").Value = 657
$ws.Range("J7").Value = 1773

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 167
$ws.Range("J7").Value = 566

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 426
$ws.Range("J6").Value = 454
$ws.Range("J7").Value = 1278

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 288
$ws.Range("J6").Value = 252
$ws.Range("J7").Value = 863

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 189
$ws.Range("J7").Value = 708

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J2").Value = 35
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 222
$ws.Range("J8").Value = 1773
$ws.Range("J11").Value = 504
$ws.Range("J14").Value = 149
$ws.Range("J15").Value = 348
$ws.Range("J18").Value = 228
$ws.Range("J19").Value = 815
$ws.Range("J20").Value = 609
$ws.Range("J23").Value = 258
$ws.Range("J24").Value = 96
$ws.Range("J25").Value = 144
$ws.Range("J29").Value = 1502
$ws.Range("J30").Value = 98
$ws.Range("J31").Value = 298
$ws.Range("J33").Value = 1278
$ws.Range("J34").Value = 131
$ws.Range("J36").Value = 383
$ws.Range("J37").Value = 863
$ws.Range("J42").Value = 1198
$ws.Range("J47").Value = 205
$ws.Range("J48").Value = 315
$ws.Range("J50").Value = 170
$ws.Range("J51").Value = 355
$ws.Range("J54").Value = 558
$ws.Range("J55").Value = 443
$ws.Range("J57").Value = 134
$ws.Range("J60").Value = 166
$ws.Range("J63").Value = 82
$ws.Range("J64").Value = 189
$ws.Range("J65").Value = 708
$ws.Range("J68").Value = 60
$ws.Range("J73").Value = 275
$ws.Range("J76").Value = 400
$ws.Range("J77").Value = 198
$ws.Range("J79").Value = 770
$ws.Range("J83").Value = 566
$ws.Range("J84").Value = 234
$ws.Range("J85").Value = 1160
$ws.Range("J89").Value = 353
$ws.Range("J91").Value = 324
$ws.Range("J92").Value = 92
$ws.Range("J94").Value = 314
$ws.Range("J96").Value = 318
$ws.Range("J100").Value = 49
$ws.Range("J101").Value = 28180

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 97
$ws.Range("J3").Value = 73
$ws.Range("J6").Value = 107
$ws.Range("J7").Value = 298

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J2").Value = 71
$ws.Range("J7").Value = 234

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 258
$ws.Range("J7").Value = 558

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 458
$ws.Range("J3").Value = 528
$ws.Range("J4").Value = 81
$ws.Range("J5").Value = 56
$ws.Range("J7").Value = 1502

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J4").Value = 49
$ws.Range("J7").Value = 315

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 200
$ws.Range("J3").Value = 231
$ws.Range("J6").Value = 316
$ws.Range("J7").Value = 815

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J3").Value = 89
$ws.Range("J7").Value = 400

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 250
$ws.Range("J3").Value = 241
$ws.Range("J6").Value = 634
$ws.Range("J7").Value = 1198

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J2").Value = 86
$ws.Range("J7").Value = 443

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J3").Value = 85
$ws.Range("J7").Value = 258

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 318

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 85
$ws.Range("J3").Value = 132
$ws.Range("J7").Value = 324

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 219
$ws.Range("J6").Value = 231
$ws.Range("J7").Value = 770

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 51
$ws.Range("J7").Value = 189

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 169
$ws.Range("J3").Value = 199
$ws.Range("J5").Value = 18
$ws.Range("J7").Value = 609

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J4").Value = 12
$ws.Range("J7").Value = 228

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 124
$ws.Range("J7").Value = 383

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 131

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 166
$ws.Range("J7").Value = 314

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 58
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J2").Value = 48
$ws.Range("J7").Value = 205

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J6").Value = 161
$ws.Range("J7").Value = 348

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J2").Value = 43
$ws.Range("J4").Value = 26
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J3").Value = 87
$ws.Range("J6").Value = 240
$ws.Range("J7").Value = 504

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 71
$ws.Range("J7").Value = 275

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J3").Value = 53
$ws.Range("J7").Value = 222

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J3").Value = 100
$ws.Range("J6").Value = 109
$ws.Range("J7").Value = 353

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J4").Value = 33
$ws.Range("J5").Value = 9
$ws.Range("J6").Value = 147
$ws.Range("J7").Value = 355

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J2").Value = 57
$ws.Range("J7").Value = 166

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 308
$ws.Range("J3").Value = 419
$ws.Range("J6").Value = 332
$ws.Range("J7").Value = 1160

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J2").Value = 74
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 201
$ws.Range("J6").Value = 307
